$wb = $excel.ActiveWorkbook

# The new "2022-Q1" sheet must slot in right before "总计", and "总计" needs to
# end up with a higher internal sheetId than "2022-Q1" (it was re-saved after the
# new sheet in the original edit). Deleting and recreating it - after inserting
# "2022-Q1" - reproduces that ordering; its handful of existing rows are simply
# re-entered (shifted down by one) along with the new first row.
$oldTotal = $wb.Worksheets.Item("总计")
$oldTotal.Delete()

# --- "2022-Q1" sheet: fund-level detail, inserted right after "2021-Q4" ---
$q4Sheet = $wb.Worksheets.Item("2021-Q4")
$newSheet = $wb.Worksheets.Add($null, $q4Sheet)
$newSheet.Name = "2022-Q1"

# Header row + index column reuse the bold/centered/bordered style already used
# by the other quarter sheets (copying formats keeps the shared style, instead of
# inventing a new one).
$q4Sheet.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)
$q4Sheet.Range("A2").Copy()
$newSheet.Range("A2:A21").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Columns B:G carry numeric-looking data (fund codes with leading zeros, fixed-
# precision percentages) that must be stored as literal text, not numbers - force
# Text format before assigning the values.
$newSheet.Range("B2:G21").NumberFormat = "@"

$newSheet.Cells.Item(2, 1).Value = 0
$newSheet.Cells.Item(2, 2).Value = '166801'
$newSheet.Cells.Item(2, 3).Value = '浙商聚潮新思维混合'
$newSheet.Cells.Item(2, 4).Value = '9.62'
$newSheet.Cells.Item(2, 5).Value = '78.06'
$newSheet.Cells.Item(2, 6).Value = '2.94'
$newSheet.Cells.Item(2, 7).Value = '0.2828'
$newSheet.Cells.Item(2, 8).Value = 9
$newSheet.Cells.Item(3, 1).Value = 1
$newSheet.Cells.Item(3, 2).Value = '009230'
$newSheet.Cells.Item(3, 3).Value = '鹏华安和混合A'
$newSheet.Cells.Item(3, 4).Value = '14.02'
$newSheet.Cells.Item(3, 5).Value = '34.45'
$newSheet.Cells.Item(3, 6).Value = '1.60'
$newSheet.Cells.Item(3, 7).Value = '0.2243'
$newSheet.Cells.Item(3, 8).Value = 3
$newSheet.Cells.Item(4, 1).Value = 2
$newSheet.Cells.Item(4, 2).Value = '009667'
$newSheet.Cells.Item(4, 3).Value = '鹏华安庆混合A'
$newSheet.Cells.Item(4, 4).Value = '11.22'
$newSheet.Cells.Item(4, 5).Value = '38.92'
$newSheet.Cells.Item(4, 6).Value = '1.86'
$newSheet.Cells.Item(4, 7).Value = '0.2087'
$newSheet.Cells.Item(4, 8).Value = 2
$newSheet.Cells.Item(5, 1).Value = 3
$newSheet.Cells.Item(5, 2).Value = '005416'
$newSheet.Cells.Item(5, 3).Value = '鹏华尊惠18个月定期开放混合A'
$newSheet.Cells.Item(5, 4).Value = '7.95'
$newSheet.Cells.Item(5, 5).Value = '37.81'
$newSheet.Cells.Item(5, 6).Value = '2.37'
$newSheet.Cells.Item(5, 7).Value = '0.1884'
$newSheet.Cells.Item(5, 8).Value = 2
$newSheet.Cells.Item(6, 1).Value = 4
$newSheet.Cells.Item(6, 2).Value = '011160'
$newSheet.Cells.Item(6, 3).Value = '富国质量成长6个月持有期混合A'
$newSheet.Cells.Item(6, 4).Value = '3.80'
$newSheet.Cells.Item(6, 5).Value = '91.55'
$newSheet.Cells.Item(6, 6).Value = '2.51'
$newSheet.Cells.Item(6, 7).Value = '0.0954'
$newSheet.Cells.Item(6, 8).Value = 4
$newSheet.Cells.Item(7, 1).Value = 5
$newSheet.Cells.Item(7, 2).Value = '009231'
$newSheet.Cells.Item(7, 3).Value = '鹏华安和混合C'
$newSheet.Cells.Item(7, 4).Value = '5.33'
$newSheet.Cells.Item(7, 5).Value = '34.45'
$newSheet.Cells.Item(7, 6).Value = '1.60'
$newSheet.Cells.Item(7, 7).Value = '0.0853'
$newSheet.Cells.Item(7, 8).Value = 3
$newSheet.Cells.Item(8, 1).Value = 6
$newSheet.Cells.Item(8, 2).Value = '006199'
$newSheet.Cells.Item(8, 3).Value = '长盛同锦研究精选混合'
$newSheet.Cells.Item(8, 4).Value = '1.73'
$newSheet.Cells.Item(8, 5).Value = '82.48'
$newSheet.Cells.Item(8, 6).Value = '4.52'
$newSheet.Cells.Item(8, 7).Value = '0.0782'
$newSheet.Cells.Item(8, 8).Value = 1
$newSheet.Cells.Item(9, 1).Value = 7
$newSheet.Cells.Item(9, 2).Value = '001892'
$newSheet.Cells.Item(9, 3).Value = '长盛新兴成长主题灵活配置混合'
$newSheet.Cells.Item(9, 4).Value = '1.32'
$newSheet.Cells.Item(9, 5).Value = '82.10'
$newSheet.Cells.Item(9, 6).Value = '4.52'
$newSheet.Cells.Item(9, 7).Value = '0.0597'
$newSheet.Cells.Item(9, 8).Value = 1
$newSheet.Cells.Item(10, 1).Value = 8
$newSheet.Cells.Item(10, 2).Value = '009242'
$newSheet.Cells.Item(10, 3).Value = '中加核心智造混合A'
$newSheet.Cells.Item(10, 4).Value = '2.05'
$newSheet.Cells.Item(10, 5).Value = '65.71'
$newSheet.Cells.Item(10, 6).Value = '2.90'
$newSheet.Cells.Item(10, 7).Value = '0.0594'
$newSheet.Cells.Item(10, 8).Value = 7
$newSheet.Cells.Item(11, 1).Value = 9
$newSheet.Cells.Item(11, 2).Value = '003165'
$newSheet.Cells.Item(11, 3).Value = '鹏华弘嘉灵活配置混合A'
$newSheet.Cells.Item(11, 4).Value = '1.53'
$newSheet.Cells.Item(11, 5).Value = '93.95'
$newSheet.Cells.Item(11, 6).Value = '3.12'
$newSheet.Cells.Item(11, 7).Value = '0.0477'
$newSheet.Cells.Item(11, 8).Value = 8
$newSheet.Cells.Item(12, 1).Value = 10
$newSheet.Cells.Item(12, 2).Value = '009668'
$newSheet.Cells.Item(12, 3).Value = '鹏华安庆混合C'
$newSheet.Cells.Item(12, 4).Value = '2.36'
$newSheet.Cells.Item(12, 5).Value = '38.92'
$newSheet.Cells.Item(12, 6).Value = '1.86'
$newSheet.Cells.Item(12, 7).Value = '0.0439'
$newSheet.Cells.Item(12, 8).Value = 2
$newSheet.Cells.Item(13, 1).Value = 11
$newSheet.Cells.Item(13, 2).Value = '002085'
$newSheet.Cells.Item(13, 3).Value = '长盛互联网+主题灵活配置混合'
$newSheet.Cells.Item(13, 4).Value = '0.84'
$newSheet.Cells.Item(13, 5).Value = '83.97'
$newSheet.Cells.Item(13, 6).Value = '4.52'
$newSheet.Cells.Item(13, 7).Value = '0.0380'
$newSheet.Cells.Item(13, 8).Value = 1
$newSheet.Cells.Item(14, 1).Value = 12
$newSheet.Cells.Item(14, 2).Value = '003166'
$newSheet.Cells.Item(14, 3).Value = '鹏华弘嘉灵活配置混合C'
$newSheet.Cells.Item(14, 4).Value = '0.56'
$newSheet.Cells.Item(14, 5).Value = '93.95'
$newSheet.Cells.Item(14, 6).Value = '3.12'
$newSheet.Cells.Item(14, 7).Value = '0.0175'
$newSheet.Cells.Item(14, 8).Value = 8
$newSheet.Cells.Item(15, 1).Value = 13
$newSheet.Cells.Item(15, 2).Value = '012005'
$newSheet.Cells.Item(15, 3).Value = '信达澳银恒盛混合A'
$newSheet.Cells.Item(15, 4).Value = '1.87'
$newSheet.Cells.Item(15, 5).Value = '31.90'
$newSheet.Cells.Item(15, 6).Value = '0.87'
$newSheet.Cells.Item(15, 7).Value = '0.0163'
$newSheet.Cells.Item(15, 8).Value = 3
$newSheet.Cells.Item(16, 1).Value = 14
$newSheet.Cells.Item(16, 2).Value = '005417'
$newSheet.Cells.Item(16, 3).Value = '鹏华尊惠18个月定期开放混合C'
$newSheet.Cells.Item(16, 4).Value = '0.56'
$newSheet.Cells.Item(16, 5).Value = '37.81'
$newSheet.Cells.Item(16, 6).Value = '2.37'
$newSheet.Cells.Item(16, 7).Value = '0.0133'
$newSheet.Cells.Item(16, 8).Value = 2
$newSheet.Cells.Item(17, 1).Value = 15
$newSheet.Cells.Item(17, 2).Value = '009128'
$newSheet.Cells.Item(17, 3).Value = '明亚价值长青混合A'
$newSheet.Cells.Item(17, 4).Value = '0.38'
$newSheet.Cells.Item(17, 5).Value = '49.48'
$newSheet.Cells.Item(17, 6).Value = '3.49'
$newSheet.Cells.Item(17, 7).Value = '0.0133'
$newSheet.Cells.Item(17, 8).Value = 5
$newSheet.Cells.Item(18, 1).Value = 16
$newSheet.Cells.Item(18, 2).Value = '011161'
$newSheet.Cells.Item(18, 3).Value = '富国质量成长6个月持有期混合C'
$newSheet.Cells.Item(18, 4).Value = '0.12'
$newSheet.Cells.Item(18, 5).Value = '91.55'
$newSheet.Cells.Item(18, 6).Value = '2.51'
$newSheet.Cells.Item(18, 7).Value = '0.0030'
$newSheet.Cells.Item(18, 8).Value = 4
$newSheet.Cells.Item(19, 1).Value = 17
$newSheet.Cells.Item(19, 2).Value = '009243'
$newSheet.Cells.Item(19, 3).Value = '中加核心智造混合C'
$newSheet.Cells.Item(19, 4).Value = '0.10'
$newSheet.Cells.Item(19, 5).Value = '65.71'
$newSheet.Cells.Item(19, 6).Value = '2.90'
$newSheet.Cells.Item(19, 7).Value = '0.0029'
$newSheet.Cells.Item(19, 8).Value = 7
$newSheet.Cells.Item(20, 1).Value = 18
$newSheet.Cells.Item(20, 2).Value = '012006'
$newSheet.Cells.Item(20, 3).Value = '信达澳银恒盛混合C'
$newSheet.Cells.Item(20, 4).Value = '0.31'
$newSheet.Cells.Item(20, 5).Value = '31.90'
$newSheet.Cells.Item(20, 6).Value = '0.87'
$newSheet.Cells.Item(20, 7).Value = '0.0027'
$newSheet.Cells.Item(20, 8).Value = 3
$newSheet.Cells.Item(21, 1).Value = 19
$newSheet.Cells.Item(21, 2).Value = '009129'
$newSheet.Cells.Item(21, 3).Value = '明亚价值长青混合C'
$newSheet.Cells.Item(21, 4).Value = '0.00'
$newSheet.Cells.Item(21, 5).Value = '49.48'
$newSheet.Cells.Item(21, 6).Value = '3.49'
$newSheet.Cells.Item(21, 7).Value = 0
$newSheet.Cells.Item(21, 8).Value = 5

# Drop the Text number format now that the literal values are stored, matching
# the default (unformatted) style the other data cells use.
$newSheet.Range("B2:G21").Style = "Normal"

# --- "总计" sheet: recreated right after "2022-Q1" ---
$newTotal = $wb.Worksheets.Add($null, $newSheet)
$newTotal.Name = "总计"

$newSheet.Range("B1:D1").Copy()
$newTotal.Range("B1:D1").PasteSpecial(-4122)
$newSheet.Range("A2").Copy()
$newTotal.Range("A1:A7").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$newTotal.Range("B1").Value = "日期"
$newTotal.Range("C1").Value = "持有数量(只)"
$newTotal.Range("D1").Value = "持有市值(亿元)"

$newTotal.Cells.Item(2, 1).Value = 0
$newTotal.Cells.Item(2, 2).Value = "2022-Q1"
$newTotal.Cells.Item(2, 3).Value = 20
$newTotal.Cells.Item(2, 4).Value = 1.48
$newTotal.Cells.Item(3, 1).Value = 1
$newTotal.Cells.Item(3, 2).Value = "2021-Q4"
$newTotal.Cells.Item(3, 3).Value = 25
$newTotal.Cells.Item(3, 4).Value = 2.23
$newTotal.Cells.Item(4, 1).Value = 2
$newTotal.Cells.Item(4, 2).Value = "2021-Q3"
$newTotal.Cells.Item(4, 3).Value = 15
$newTotal.Cells.Item(4, 4).Value = 2.05
$newTotal.Cells.Item(5, 1).Value = 3
$newTotal.Cells.Item(5, 2).Value = "2021-Q2"
$newTotal.Cells.Item(5, 3).Value = 9
$newTotal.Cells.Item(5, 4).Value = 1.2
$newTotal.Cells.Item(6, 1).Value = 4
$newTotal.Cells.Item(6, 2).Value = "2021-Q1"
$newTotal.Cells.Item(6, 3).Value = 15
$newTotal.Cells.Item(6, 4).Value = 5.88
$newTotal.Cells.Item(7, 1).Value = 5
$newTotal.Cells.Item(7, 2).Value = "2020-Q4"
$newTotal.Cells.Item(7, 3).Value = 7
$newTotal.Cells.Item(7, 4).Value = 2.56

Write-Output "edit complete"
